$wb = $excel.ActiveWorkbook

# --- Users sheet (sheet1.xml): add an "Id" column before Username/Password ---
$wsUsers = $wb.Worksheets.Item("Users")

# Insert a new column A, shifting existing Username/Password data to B/C
$wsUsers.Columns.Item(1).Insert()

# Write the row values before the header, so the shared-string table picks up
# "u01"/"u02"/"u03" ahead of "Id" (matches the target string order)
$wsUsers.Range("A2").Value = "u01"
$wsUsers.Range("A3").Value = "u02"
$wsUsers.Range("A4").Value = "u03"
$wsUsers.Range("A1").Value = "Id"

# Make the Users sheet the active tab/selected sheet, set its zoom and selection
$wsUsers.Activate()
$wsUsers.Range("D6").Select()
$excel.ActiveWindow.Zoom = 190

$wb.Save()
